$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSE")

# Update the RSE sheet values per the feed correction
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 10
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 1000

# Make RSE the active sheet/tab and select C3 (moves tabSelected + active
# cell from R2M -> RSE, and bumps the workbook's activeTab)
$ws.Activate()
$ws.Range("C3").Select()
